$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (id) and C (speaker_variant) for rows 2-17
$values = @(
    @(2,  "#rabbi",        "Rabbi"),
    @(3,  "#gelasius",     "Gelasius"),
    @(4,  "#lechi",        "Lechi"),
    @(5,  "#president",    "President"),
    @(6,  "#iason",        "Iason"),
    @(7,  "#eurialus",     "Eurialus"),
    @(8,  "#cleanthes",    "Cleanthes"),
    @(9,  "#titus",        "Titus"),
    @(10, "#epistemon",    "Epistemon"),
    @(11, "#de-president", "De President"),
    @(12, "#paulus",       "Paulus"),
    @(13, "#epicurus",     "Epicurus"),
    @(14, "#lechi",        "lechi"),
    @(15, "#choor",        "Choor"),
    @(16, "#talus",        "Talus"),
    @(17, "#corydon",      "Corydon")
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Clear the "is_prefered" (column D) marker for all data rows, since the
# export no longer marks a preferred variant.
$ws.Range("D2:D15").ClearContents()
